# LSTM_mv.xlsx update — append latest day rows of market data to each
# tracking sheet (D1_USD, D1_EUR, D5_EUR, D1_OIL), fix the D5_EUR "Nan"
# placeholders for days that now have predictions, and leave the
# workbook with D5_EUR as the active/selected sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# D1_USD  (sheet1) — add rows 130:133, correct B129
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("D1_USD")

# Row 129's predicted value was revised.
$ws.Range("B129").Value = 4.021039

# Clone formatting from the last existing data row down across the new ones.
$ws.Range("A129:E129").Copy()
$ws.Range("A130:E132").PasteSpecial(-4122)

$ws.Range("A130").Value = 45342
$ws.Range("B130").Value = 4.0136430000000001
$ws.Range("C130").Value = 4.027825

$ws.Range("A131").Value = 45343
$ws.Range("B131").Value = 3.9898500000000001
$ws.Range("C131").Value = 4.0217590000000003

$ws.Range("D130:D131").Formula = "=B130-C130"
$ws.Range("E130:E131").Formula = "=IF(D130<0,1,0)"

$ws.Range("A132").Value = 45344
$ws.Range("B132").Value = 3.9882620000000002
$ws.Range("C132").Value = 3.9474418
$ws.Range("D132").Formula = "=B132-C132"
$ws.Range("E132").Formula = "=IF(D132<0,1,0)"

# Row 133 only carries a stray C value (next day's actual, no prediction yet).
$ws.Range("C133").Value = 4.0039740000000004

# ---------------------------------------------------------------------
# D1_EUR  (sheet3) — add rows 405:407, fill out row 404
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("D1_EUR")

$ws.Range("A403:E403").Copy()
$ws.Range("A404:E406").PasteSpecial(-4122)

$ws.Range("A404").Value = 45342
$ws.Range("B404").Value = 4.3247
$ws.Range("D404").Formula = "=B404-C404"
$ws.Range("E404").Formula = "=IF(D404<0,1,0)"

$ws.Range("A405").Value = 45343
$ws.Range("B405").Value = 4.3125299999999998
$ws.Range("C405").Value = 4.3379659999999998
$ws.Range("D405").Formula = "=B405-C405"
$ws.Range("E405").Formula = "=IF(D405<0,1,0)"

$ws.Range("A406").Value = 45344
$ws.Range("B406").Value = 4.3152999999999997
$ws.Range("C406").Value = 4.3258729999999996
$ws.Range("D406").Formula = "=B406-C406"
$ws.Range("E406").Formula = "=IF(D406<0,1,0)"

$ws.Range("C407").Value = 4.3224660000000004

# ---------------------------------------------------------------------
# D5_EUR  (sheet5) — rows 96:99 predictions resolved from "Nan" to values,
# add rows 100:104 with fresh "Nan" placeholders
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("D5_EUR")

$ws.Range("B96").Value = 4.3365299999999998
$ws.Range("B97").Value = 4.3247
$ws.Range("B98").Value = 4.3125299999999998
$ws.Range("B99").Value = 4.3152999999999997

$ws.Range("A99:C99").Copy()
$ws.Range("A100:C104").PasteSpecial(-4122)

$ws.Range("A100").Value = 45345
$ws.Range("B100").Value = "Nan"
$ws.Range("C100").Value = 4.3074820000000003

$ws.Range("A101").Value = 45348
$ws.Range("B101").Value = "Nan"
$ws.Range("C101").Value = 4.3080907000000002

$ws.Range("A102").Value = 45349
$ws.Range("B102").Value = "Nan"
$ws.Range("C102").Value = 4.3042490000000004

$ws.Range("A103").Value = 45350
$ws.Range("B103").Value = "Nan"
$ws.Range("C103").Value = 4.3033146999999996

$ws.Range("A104").Value = 45351
$ws.Range("B104").Value = "Nan"
$ws.Range("C104").Value = 4.3039784000000001

# ---------------------------------------------------------------------
# D1_OIL  (sheet6) — add rows 77:79, fill out row 76
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("D1_OIL")

$ws.Range("A75:E75").Copy()
$ws.Range("A76:E78").PasteSpecial(-4122)

$ws.Range("A76").Value = 45342
$ws.Range("B76").Value = 78.269997000000004
$ws.Range("D76").Formula = "=B76-C76"
$ws.Range("E76").Formula = "=D76/C76"

$ws.Range("A77").Value = 45343
$ws.Range("B77").Value = 77.910004000000001
$ws.Range("C77").Value = 67.810500000000005
$ws.Range("D77").Formula = "=B77-C77"
$ws.Range("E77").Formula = "=D77/C77"

$ws.Range("A78").Value = 45344
$ws.Range("B78").Value = 78.610000999999997
$ws.Range("C78").Value = 68.6648
$ws.Range("D78").Formula = "=B78-C78"
$ws.Range("E78").Formula = "=D78/C78"

$ws.Range("C79").Value = 69.136099999999999

# ---------------------------------------------------------------------
# Final view state: D5_EUR is the sheet left on-screen/selected.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("D5_EUR")
$ws.Activate()
$ws.Range("B100").Select()
